# Add the new user "Juanita" (role VENDEDOR) as row 5 of the Datos sheet,
# right after the existing rows (usuario / clave / rol columns A:C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Juanita"

# "clave" column stores passwords as text (the existing rows store "1234"
# as text too), so force text with a leading apostrophe instead of letting
# Excel coerce it to a number, then restore the default "Normal" style so
# the new cell matches the unstyled look of the rest of the table.
$ws.Range("B5").Value = "'1234"
$ws.Range("B5").Style = "Normal"

$ws.Range("C5").Value = "VENDEDOR"

# Keep Excel's "number stored as text" warning suppressed over the table,
# now that it spans through row 5.
$ws.Range("A1:C5").Errors.Item(9).Ignore = $true
